$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 243.58333
$ws.Range("I2").Value = 259.2
$ws.Range("J2").Value = 232.42857
$ws.Range("K2").Value = 259.2
$ws.Range("L2").Value = 232.42857
$ws.Range("M2").Value = -146.2
$ws.Range("N2").Value = -458.42857

$ws.Range("H18").Value = 13619.111
$ws.Range("I18").Value = 5266
$ws.Range("J18").Value = 17795.666
$ws.Range("K18").Value = 5266
$ws.Range("L18").Value = 17795.666
$ws.Range("M18").Value = -4982
$ws.Range("N18").Value = -18363.666

$ws.Range("H33").Value = 95.92308
$ws.Range("I33").Value = 75.5
$ws.Range("K33").Value = 75.5
$ws.Range("M33").Value = 153.5

$ws.Range("H99").Value = 941.3333
$ws.Range("I99").Value = 281.7143
$ws.Range("K99").Value = 845.1428999999999
$ws.Range("M99").Value = 652.8571000000001

$ws.Range("H101").Value = 1778987.4
$ws.Range("I101").Value = 4444594
$ws.Range("J101").Value = 1916.3334
$ws.Range("K101").Value = 13333782
$ws.Range("L101").Value = 5749.0002
$ws.Range("M101").Value = -13332160
$ws.Range("N101").Value = -8993.0002

$ws.Range("H129").Value = 876.475
$ws.Range("I129").Value = 722.25
$ws.Range("J129").Value = 893.6111
$ws.Range("K129").Value = 2166.75
$ws.Range("L129").Value = 2680.8333
$ws.Range("M129").Value = 2833.25
$ws.Range("N129").Value = -12680.8333

$ws.Range("H132").Value = 1221.02
$ws.Range("I132").Value = 1145.6383
$ws.Range("K132").Value = 3436.9149
$ws.Range("M132").Value = -906.9149000000002

$ws.Range("H137").Value = 1001.74194
$ws.Range("I137").Value = 794.9655
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 2384.8965
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = 165.1035000000002
$ws.Range("N137").Value = -17100

$ws.Range("H138").Value = 4071.3635
$ws.Range("I138").Value = 4320.7856
$ws.Range("J138").Value = 3634.875
$ws.Range("K138").Value = 12962.3568
$ws.Range("L138").Value = 10904.625
$ws.Range("M138").Value = -7822.356800000001
$ws.Range("N138").Value = -21184.625

$ws.Range("H141").Value = 1557931.4
$ws.Range("I141").Value = 2334491.5
$ws.Range("J141").Value = 4811
$ws.Range("K141").Value = 7003474.5
$ws.Range("L141").Value = 14433
$ws.Range("M141").Value = -6998294.5
$ws.Range("N141").Value = -24793

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 52003
$ws.Range("I23").Value = 80006
$ws.Range("J23").Value = 38001.5
$ws.Range("K23").Value = 80006
$ws.Range("L23").Value = 38001.5
$ws.Range("M23").Value = -79747
$ws.Range("N23").Value = -38519.5

$ws.Range("H32").Value = 3379.8333
$ws.Range("I32").Value = 2918.738
$ws.Range("K32").Value = 2918.738
$ws.Range("M32").Value = -2631.738

$ws.Range("H37").Value = 4000
$ws.Range("I37").Value = 4000
$ws.Range("K37").Value = 4000
$ws.Range("M37").Value = -3727

$ws.Range("H45").Value = 1441.6154
$ws.Range("I45").Value = 1212.7333
$ws.Range("K45").Value = 1212.7333
$ws.Range("M45").Value = -835.7333000000001

$ws.Range("H97").Value = 999.1539
$ws.Range("I97").Value = 999.1539
$ws.Range("K97").Value = 999.1539
$ws.Range("M97").Value = -503.1539

$ws.Range("H110").Value = 2342.5908
$ws.Range("J110").Value = 3220.8572
$ws.Range("L110").Value = 3220.8572
$ws.Range("N110").Value = -7310.8572

$ws.Range("H122").Value = 2064.8
$ws.Range("I122").Value = 2064.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6194.400000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3744.400000000001
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 1828.8572
$ws.Range("I132").Value = 1247.409
$ws.Range("J132").Value = 2812.8462
$ws.Range("K132").Value = 3742.227
$ws.Range("L132").Value = 8438.5386
$ws.Range("M132").Value = -1212.227
$ws.Range("N132").Value = -13498.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 7332.6665
$ws.Range("I12").Value = 6999
$ws.Range("J12").Value = 8000
$ws.Range("K12").Value = 6999
$ws.Range("L12").Value = 8000
$ws.Range("M12").Value = -6831
$ws.Range("N12").Value = -8336

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H105").Value = 2316.5
$ws.Range("I105").Value = 2218.9048
$ws.Range("K105").Value = 2218.9048
$ws.Range("M105").Value = -471.9047999999998

$ws.Range("H126").Value = 40000
$ws.Range("J126").Value = 40000
$ws.Range("L126").Value = 40000
$ws.Range("N126").Value = -49880

$ws.Range("H134").Value = 3779.9246
$ws.Range("I134").Value = 4415.5854
$ws.Range("J134").Value = 1608.0834
$ws.Range("K134").Value = 13246.7562
$ws.Range("L134").Value = 4824.2502
$ws.Range("M134").Value = -10711.7562
$ws.Range("N134").Value = -9894.2502

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1829.6111
$ws.Range("J31").Value = 2399.25
$ws.Range("L31").Value = 2399.25
$ws.Range("N31").Value = -2989.25

$ws.Range("H34").Value = 1829.6111
$ws.Range("J34").Value = 2399.25
$ws.Range("L34").Value = 2399.25
$ws.Range("N34").Value = -2803.25

$ws.Range("H58").Value = 1088093.4
$ws.Range("I58").Value = 1553712
$ws.Range("J58").Value = 1650
$ws.Range("K58").Value = 1553712
$ws.Range("L58").Value = 1650
$ws.Range("M58").Value = -1553509
$ws.Range("N58").Value = -2056

$ws.Range("H122").Value = 3798.5454
$ws.Range("I122").Value = 2122.6667
$ws.Range("K122").Value = 6368.000100000001
$ws.Range("M122").Value = -3918.000100000001

$ws.Range("H132").Value = 1477.3448
$ws.Range("I132").Value = 978.1667
$ws.Range("J132").Value = 3873.4
$ws.Range("K132").Value = 2934.5001
$ws.Range("L132").Value = 11620.2
$ws.Range("M132").Value = -404.5001000000002
$ws.Range("N132").Value = -16680.2

$ws.Range("H134").Value = 1925.4572
$ws.Range("I134").Value = 1743.08
$ws.Range("K134").Value = 5229.24
$ws.Range("M134").Value = -2694.24

$ws.Range("H136").Value = 1088093.4
$ws.Range("I136").Value = 1553712
$ws.Range("J136").Value = 1650
$ws.Range("K136").Value = 4661136
$ws.Range("L136").Value = 4950
$ws.Range("M136").Value = -4658586
$ws.Range("N136").Value = -10050

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 850.05
$ws.Range("J122").Value = 997.2857
$ws.Range("L122").Value = 8975.5713
$ws.Range("N122").Value = -13875.5713

$ws.Range("H131").Value = 13553.921
$ws.Range("J131").Value = 14896.263
$ws.Range("L131").Value = 44688.789
$ws.Range("N131").Value = -54768.789

$ws.Range("H132").Value = 886.4286
$ws.Range("J132").Value = 1200.3334
$ws.Range("L132").Value = 10803.0006
$ws.Range("N132").Value = -15863.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 13869.25
$ws.Range("J123").Value = 13869.25
$ws.Range("L123").Value = 13869.25
$ws.Range("N123").Value = -18769.25

$ws.Range("H132").Value = 1014447.8
$ws.Range("I132").Value = 1604182.2
$ws.Range("J132").Value = 3474.4285
$ws.Range("K132").Value = 4812546.6
$ws.Range("L132").Value = 10423.2855
$ws.Range("M132").Value = -4810016.6
$ws.Range("N132").Value = -15483.2855

$ws.Range("H133").Value = 85000
$ws.Range("J133").Value = 85000
$ws.Range("L133").Value = 85000
$ws.Range("N133").Value = -95120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 513.3333
$ws.Range("I55").Value = 397
$ws.Range("J55").Value = 746
$ws.Range("K55").Value = 397
$ws.Range("L55").Value = 746
$ws.Range("M55").Value = -224
$ws.Range("N55").Value = -1092

$ws.Range("H93").Value = 951.36365
$ws.Range("I93").Value = 806.5
$ws.Range("J93").Value = 2400
$ws.Range("K93").Value = 806.5
$ws.Range("L93").Value = 2400
$ws.Range("M93").Value = 441.5
$ws.Range("N93").Value = -4896

$ws.Range("H122").Value = 4183.25
$ws.Range("I122").Value = 3532.8462
$ws.Range("J122").Value = 7001.6665
$ws.Range("K122").Value = 10598.5386
$ws.Range("L122").Value = 21004.9995
$ws.Range("M122").Value = -8148.5386
$ws.Range("N122").Value = -25904.9995

$ws.Range("H132").Value = 1514.4054
$ws.Range("I132").Value = 1019.21875
$ws.Range("J132").Value = 4683.6
$ws.Range("K132").Value = 3057.65625
$ws.Range("L132").Value = 14050.8
$ws.Range("M132").Value = -527.65625
$ws.Range("N132").Value = -19110.8

$ws.Range("H136").Value = 2267.3333
$ws.Range("I136").Value = 1401.9286
$ws.Range("K136").Value = 4205.7858
$ws.Range("M136").Value = -1655.7858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 65932.086
$ws.Range("I122").Value = 71816.82000000001
$ws.Range("K122").Value = 215450.46
$ws.Range("M122").Value = -213000.46

$ws.Range("H132").Value = 1565.4722
$ws.Range("J132").Value = 2664.3076
$ws.Range("L132").Value = 7992.9228
$ws.Range("N132").Value = -13052.9228
